$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.694.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.49%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.699.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'673.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.78%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'161.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.42%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.03%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.498"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.02%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.89%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'7.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.75%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.444"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +2.16%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +1.00%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'32.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.92%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'3.694.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.39%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'69.673.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.46%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +1.72%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'16.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.47%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  +1.88%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'474.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.09%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -1.13%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D22").Value = "'80.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.97%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'3.847.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.73%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "'PEPE"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'0.0000127"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +3.45%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("B25").Value = "'Dai"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.05%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.83%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.46%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E29").Value = "'  +0.94%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +1.99%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +0.48%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +0.04%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'26.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.32%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.166"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +3.70%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'3.689.16"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'8.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +4.67%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'6.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.04%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +0.01%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -0.05%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.52%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.0907"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.17%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'172.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +3.80%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.940"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.04%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'47.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.92%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +2.76%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.000279"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.38%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'27.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.88%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -0.68%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -0.19%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +1.79%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +0.77%  "
$ws.Range("E51").Style = "Normal"
